# correct small and medium mesh size differences between aug 2019 and jan 2020
#
# raw_data_1 had the sm_rf_kg (column L, "5x5cm seive") and the new
# md_rf_kg (column K, "2.5x2.5cm seive") values entered in swapped columns
# for the rows collected in Aug 2019 (rows 2-65). This script swaps the
# K/L values back into their correct columns, fixes the corresponding
# cross-check formula on raw_check, and adds the missing md_rf_kg
# description row (and its shared-string text) on the description sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. raw_data_1: swap columns K (md_rf_kg) and L (sm_rf_kg) for rows 2-65
# ---------------------------------------------------------------------
$rawData1 = $wb.Worksheets.Item("raw_data_1")

for ($r = 2; $r -le 65; $r++) {
    $kVal = $rawData1.Cells.Item($r, 11).Value2
    $lVal = $rawData1.Cells.Item($r, 12).Value2
    $rawData1.Cells.Item($r, 11).Value2 = $lVal
    $rawData1.Cells.Item($r, 12).Value2 = $kVal
}

# ---------------------------------------------------------------------
# 2. raw_check: the K-column cross check formula compared raw_data_1's L
#    column against raw_data_2's L column; it should compare raw_data_1's
#    (now corrected) K column against raw_data_2's L column.
# ---------------------------------------------------------------------
$rawCheck = $wb.Worksheets.Item("raw_check")

for ($r = 2; $r -le 65; $r++) {
    $rawCheck.Cells.Item($r, 11).Formula = "=IF(raw_data_1!K$r=raw_data_2!L$r,`"`",`"check`")"
}

# ---------------------------------------------------------------------
# 3. description: re-point the existing sm_rf_kg description row at
#    md_rf_kg, then insert a new row directly below it with the real
#    md_rf_kg description (collector/recorder rows shift down by one).
# ---------------------------------------------------------------------
$description = $wb.Worksheets.Item("description")

$description.Rows.Item(12).Insert()
$description.Range("A11").Value2 = "md_rf_kg"
$description.Range("A12").Value2 = "md_rf_kg"
$description.Range("B12").Value2 = "mass of the weighing bucket plus all reef material retained in the 2.5 x 2.5cm seive in kg.  This excludes the material retained in the 5cm and 10cm seives."

# update the selection left on the description sheet
$description.Activate()
$description.Range("G13").Select()

# ---------------------------------------------------------------------
# 4. leave raw_data_1 as the active sheet/selection, scrolled further
#    down and with the cursor back at the top of the sheet.
# ---------------------------------------------------------------------
$rawData1.Activate()
$rawData1.Range("A99").Select()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
